# Adds 3 more "test" rows (rows 10, 11 and 12) worth of timing data to the
# ASSESSOR+ timesheet, plus fills in the two previously-empty rows (23-25)
# further down the sheet with their start/stop recording times.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10 : Feature 8 -----------------------------------------------
$ws.Range("F10:G10").NumberFormat = "h:mm"
$ws.Range("F10").Value = 0.35694444444444445
$ws.Range("G10").Value = 0.3576388888888889
$ws.Range("H10").Value = "Correct"
$ws.Range("I10").Value = "Correct"
$ws.Range("J10").Value = "Correct: The generate code was correct."

# --- Row 11 : Feature 9 -----------------------------------------------
$ws.Range("F11:G11").NumberFormat = "h:mm"
$ws.Range("F11").Value = 0.36388888888888887
$ws.Range("G11").Value = 0.3659722222222222
$ws.Range("H11").Value = "Correct"
$ws.Range("I11").Value = "Correct"
$ws.Range("J11").Value = "Correct: The generate code was correct."

# --- Row 12 : Feature 10 (bottom, thick-bordered row) ------------------
$ws.Range("F12:G12").NumberFormat = "h:mm"
$ws.Range("F12").Value = 0.37152777777777779
$ws.Range("G12").Value = 0.37430555555555556
$ws.Range("H12").Value = "Correct"
$ws.Range("I12").Value = "Correct"
$ws.Range("J12").Value = "Correct: The generate code was correct."

# --- Second table (Prestashop, MANUAL rows) : rows 23-25 ---------------
$ws.Range("F23:G23").NumberFormat = "h:mm"
$ws.Range("F23").Value = 0.62916666666666665
$ws.Range("G23").Value = 0.63749999999999996

$ws.Range("F24:G25").NumberFormat = "h:mm"
$ws.Range("F24").Value = 0.64930555555555558
$ws.Range("G24").Value = 0.65763888888888888
$ws.Range("F25").Value = 0.67777777777777781
$ws.Range("G25").Value = 0.6791666666666667

# Move/restore the active selection to where the author left off editing.
$ws.Range("J19").Select()
